# Automatische test-sync: 2025-07-27 19:52:50
# Adds a new test-mail log entry (row 21) to the "Logs" sheet, extends the
# conditional formatting ranges that covered rows 2-20 to also cover row 21,
# and updates the "Dashboard" category-count table to reflect the fact that
# "Intern verzoek / Actie voor medewerker" now ties with "Productinformatie".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# --- New row 21 on the "Logs" sheet --------------------------------------
$ws.Range("A21").Value = "Bel jij klant Jansen even?"
$ws.Range("B21").Value = "mailmind.test@zohomail.eu"
$ws.Range("C21").Value = "Testmail #19: Bel jij klant Jansen even?"
$ws.Range("D21").Value = "Intern verzoek / Actie voor medewerker"
$ws.Range("E21").Value = "Beste [naam],`nDank voor je e-mail. Om je verzoek efficiënt te kunnen verwerken, zouden we graag wat meer informatie ontvangen. Zou je zo vriendelijk willen zijn om de contactgegevens van klant Jansen met ons te delen, zodat we contact met hem kunnen opnemen?`nMet vriendelijke groet,`n[Bedrijfsnaam]"
$ws.Range("F21").Value = "2025-07-27 19:52:27"
$ws.Range("G21").Value = "Ja"
$ws.Range("H21").Value = "Nee"
$ws.Range("I21").Value = "Ja"
$ws.Range("J21").Value = "Nee"

# Multi-line content in E21 makes Excel pin an explicit row height
# (customHeight="1"); AutoFit clears that back to the sheet default so the
# row matches the plain, unattributed <row r="21"> the source file has.
$ws.Rows.Item(21).AutoFit()

# --- Extend conditional formatting sqref from row 20 to row 21 ----------
# Every cfRule that shares the same "applies to" range gets moved at once,
# so touching one rule per column-group is enough.
$ws.Range("D2:D20").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D21"))
$ws.Range("G2:G20").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G21"))
$ws.Range("H2:H20").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H2:H21"))
$ws.Range("I2:I20").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("I2:I21"))
$ws.Range("J2:J20").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("J2:J21"))

# --- Update the "Dashboard" category-count table -------------------------
# "Intern verzoek / Actie voor medewerker" count rises from 4 to 5 (because
# of the new row above) and now ties with "Productinformatie", so it moves
# up into row 3 while "Productinformatie" drops to row 4.
$ws2 = $wb.Worksheets.Item("Dashboard")
$ws2.Range("A3").Value = "Intern verzoek / Actie voor medewerker"
$ws2.Range("B3").Value = 5
$ws2.Range("A4").Value = "Productinformatie"
$ws2.Range("B4").Value = 5
